# Applies the RoostGPT-generated functional-test content to test-new-hello.xlsx:
#  1. Renames the worksheet from "Test Scenarios" to "Sheet1"
#  2. Widens column B to match the other (50-wide) columns
#  3. Appends six new test-case rows (TC-012..TC-017) covering the new
#     international-remittance scenarios

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet
$ws.Name = "Sheet1"

# 2) Column B was a narrow 15-wide "testId" column; bring it in line with the
#    rest (50 wide). 49.14 is the ColumnWidth value Excel persists as width="50".
$ws.Columns.Item(2).ColumnWidth = 49.14

# 3) New test-case rows, one array of six cell strings per row (A..F):
#    scenario, testId, testDescription, prerequisites, stepsToPerform, expectedResult
$newRows = @()

# Row 13 - TC-012
$row = @(
@'
<Tier 2 User Successfully Sends International Remittance-Verify a Tier 2 (Verified) user can successfully initiate an international money transfer.>
'@,
@'
TC-012
'@,
@'
This test case validates the happy path for the international remittance feature as per US-401.
'@,
@'
User must be logged in as a 'Tier 2 (Verified)' user with sufficient funds in their account. A beneficiary must be set up.
'@,
@'
1. Log in as a Tier 2 user.
2. Navigate to the 'International Transfer' or 'Send Money' section.
3. Select a pre-saved beneficiary.
4. Enter the amount to send.
5. Review the exchange rate and fees.
6. Confirm the transaction.
'@,
@'
The transfer is initiated successfully. The user sees a confirmation screen with a transaction reference number. The account balance is updated, and the transfer appears in the international transfer history with a 'Pending' status.
'@
)
$newRows += ,$row

# Row 14 - TC-013
$row = @(
@'
<Tier 1 User Blocked from International Remittance-Verify a Tier 1 (Unverified) user is not able to access or use the international remittance feature.>
'@,
@'
TC-013
'@,
@'
This negative test case enforces the business rule in US-401 that only verified users can send money internationally.
'@,
@'
User must be logged in as a 'Tier 1 (Unverified)' user.
'@,
@'
1. Log in as a Tier 1 user.
2. Attempt to navigate to the 'International Transfer' or 'Send Money' section.
'@,
@'
The 'International Transfer' option should be disabled or not visible. If accessed directly, the user should be redirected or shown a message stating they need to be a Tier 2 user to access this feature.
'@
)
$newRows += ,$row

# Row 15 - TC-014
$row = @(
@'
<User Adds and Saves a New Beneficiary-Verify a user can add and save the details of a new international beneficiary.>
'@,
@'
TC-014
'@,
@'
This test case validates the beneficiary management feature described in US-402.
'@,
@'
User must be logged in as a Tier 2 user.
'@,
@'
1. Log in as a Tier 2 user.
2. Navigate to the 'Beneficiaries' or 'Recipients' management page.
3. Click 'Add New Beneficiary'.
4. Fill in all required fields (Full Name, Country, Bank Account/Mobile Money details).
5. Save the beneficiary.
'@,
@'
The new beneficiary is saved successfully and appears in the list of saved beneficiaries. A confirmation message is displayed.
'@
)
$newRows += ,$row

# Row 16 - TC-015
$row = @(
@'
<User Deletes an Existing Beneficiary-Verify a user can delete a previously saved beneficiary.>
'@,
@'
TC-015
'@,
@'
This test case validates the beneficiary deletion functionality from US-402.
'@,
@'
User must be logged in as a Tier 2 user and have at least one saved beneficiary.
'@,
@'
1. Log in as a Tier 2 user.
2. Navigate to the 'Beneficiaries' management page.
3. Select a beneficiary from the list.
4. Click the 'Delete' or 'Remove' option.
5. Confirm the deletion in the confirmation prompt.
'@,
@'
The beneficiary is successfully removed from the list of saved beneficiaries.
'@
)
$newRows += ,$row

# Row 17 - TC-016
$row = @(
@'
<Verify Display of Exchange Rate and Fees Before Transfer Confirmation-Verify that the exchange rate and all applicable fees are clearly displayed to the user before they confirm an international transfer.>
'@,
@'
TC-016
'@,
@'
This test case validates the transparency requirement of US-403.
'@,
@'
User is logged in as a Tier 2 user and is in the process of initiating an international transfer.
'@,
@'
1. Log in as a Tier 2 user.
2. Navigate to the 'International Transfer' section.
3. Select a beneficiary and enter a sending amount.
4. Proceed to the confirmation/review screen.
'@,
@'
The confirmation screen must clearly display the send amount, the exchange rate being used, the calculated recipient amount, a breakdown of all transaction fees, and the total amount to be debited. The 'Confirm' button should be present.
'@
)
$newRows += ,$row

# Row 18 - TC-017
$row = @(
@'
<Verify International Transfer History and Pagination-Verify the user can view a paginated history of their international transfers with correct statuses.>
'@,
@'
TC-017
'@,
@'
This test case validates the requirements of US-404 for viewing transfer history.
'@,
@'
User must be logged in and have a history of more than 25 international transfers with various statuses (Completed, Failed, Pending).
'@,
@'
1. Log in to the platform.
2. Navigate to the 'International Transfer History' page.
3. Observe the list of transactions and their statuses.
4. Verify the number of items on the first page.
5. Click the 'Next' page button.
'@,
@'
The page displays a list of past international transfers with their status. The first page shows a maximum of 25 transactions. Clicking 'Next' loads the subsequent set of transactions.
'@
)
$newRows += ,$row

$startRow = 13
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $colLetter = [char](65 + $c)
        $ws.Range("$colLetter$r").Value = $row[$c]
    }
}
